# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the default "Table_0" style ({9EC34E4F-4582-4561-A9D2-D0569343EB1D})
#    to {39C30118-D0FA-43F4-905C-032B74D21A1F}.
# 2) The presentation's colour theme is switched from the "Red Violet"
#    (Integral) palette to the stock "Office" palette - i.e. every slot in
#    the active ThemeColorScheme is repainted with the Office theme's RGB
#    values.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-point the three tables at the new table style.
# ---------------------------------------------------------------------------
$newTableStyleId = "{39C30118-D0FA-43F4-905C-032B74D21A1F}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the active theme's colour scheme over to the "Office" palette.
# ---------------------------------------------------------------------------
function ConvertHexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the OOXML <a:clrScheme> slot order, and therefore the
# ThemeColorScheme.Colors(1..12) index order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertHexToRgbInt $officeThemeColors[$i - 1]
}
